$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6448121666908264
$ws.Range("B1").Value = 1.478277683258057
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.305998802185059
$ws.Range("E1").Value = 1.371301531791687
